$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then new values for columns B/C/D/E (only the
# columns that actually changed are present, others are $null/skipped).
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = "63.050.22"; E = "  +0.69%  " }
    @{ Row = 3; B = $null; C = $null; D = "2.445.32"; E = "  +0.37%  " }
    @{ Row = 4; B = $null; C = $null; D = $null; E = "  -0.08%  " }
    @{ Row = 5; B = $null; C = $null; D = "570.76"; E = "  +0.69%  " }
    @{ Row = 6; B = $null; C = $null; D = "145.92"; E = "  +0.42%  " }
    @{ Row = 7; B = $null; C = $null; D = $null; E = "  +0.07%  " }
    @{ Row = 8; B = $null; C = $null; D = $null; E = "  +0.66%  " }
    @{ Row = 9; B = $null; C = $null; D = "2.441.97"; E = "  +0.19%  " }
    @{ Row = 10; B = $null; C = $null; D = $null; E = "  +0.06%  " }
    @{ Row = 11; B = $null; C = $null; D = $null; E = "  +1.21%  " }
    @{ Row = 12; B = $null; C = $null; D = $null; E = "  -0.73%  " }
    @{ Row = 13; B = $null; C = $null; D = $null; E = "  +0.04%  " }
    @{ Row = 14; B = $null; C = $null; D = "26.89"; E = "  +0.29%  " }
    @{ Row = 15; B = $null; C = $null; D = $null; E = "  -1.18%  " }
    @{ Row = 16; B = $null; C = $null; D = $null; E = "  +0.10%  " }
    @{ Row = 17; B = $null; C = $null; D = "63.139.97"; E = "  +1.13%  " }
    @{ Row = 18; B = $null; C = $null; D = "2.449.85"; E = "  +0.48%  " }
    @{ Row = 19; B = $null; C = $null; D = "11.28"; E = "  +0.34%  " }
    @{ Row = 20; B = $null; C = $null; D = "7.30"; E = "  +4.86%  " }
    @{ Row = 21; B = $null; C = $null; D = "327.28"; E = "  +0.95%  " }
    @{ Row = 22; B = $null; C = $null; D = $null; E = "  +0.42%  " }
    @{ Row = 23; B = $null; C = $null; D = $null; E = "  +12.14%  " }
    @{ Row = 24; B = $null; C = $null; D = $null; E = "  +0.13%  " }
    @{ Row = 25; B = $null; C = $null; D = "65.61"; E = "  -2.47%  " }
    @{ Row = 26; B = $null; C = $null; D = "610.95"; E = "  +5.50%  " }
    @{ Row = 27; B = $null; C = $null; D = "8.87"; E = "  +3.58%  " }
    @{ Row = 28; B = $null; C = $null; D = "0.0000101"; E = "  +1.03%  " }
    @{ Row = 29; B = $null; C = $null; D = "2.577.78"; E = $null }
    @{ Row = 30; B = $null; C = $null; D = $null; E = "  +0.33%  " }
    @{ Row = 31; B = $null; C = $null; D = $null; E = "  +3.43%  " }
    @{ Row = 32; B = $null; C = $null; D = "8.20"; E = "  -2.58%  " }
    @{ Row = 33; B = $null; C = $null; D = $null; E = "  -3.42%  " }
    @{ Row = 34; B = $null; C = $null; D = $null; E = "  +0.82%  " }
    @{ Row = 35; B = $null; C = $null; D = "5.17"; E = "  +6.76%  " }
    @{ Row = 36; B = $null; C = $null; D = $null; E = "  +0.77%  " }
    @{ Row = 37; B = $null; C = $null; D = $null; E = "  +0.13%  " }
    @{ Row = 38; B = $null; C = $null; D = "0.378"; E = "  -0.93%  " }
    @{ Row = 39; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "18.74"; E = "  -0.12%  " }
    @{ Row = 40; B = "RenderToken"; C = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D = "5.38"; E = "  +0.07%  " }
    @{ Row = 41; B = $null; C = $null; D = "145.29"; E = "  -1.58%  " }
    @{ Row = 42; B = $null; C = $null; D = $null; E = "  -2.13%  " }
    @{ Row = 43; B = $null; C = $null; D = "2.58"; E = "  +6.01%  " }
    @{ Row = 44; B = $null; C = $null; D = $null; E = "  -0.08%  " }
    @{ Row = 45; B = $null; C = $null; D = "41.87"; E = "  +0.62%  " }
    @{ Row = 46; B = $null; C = $null; D = "148.25"; E = "  +0.11%  " }
    @{ Row = 47; B = $null; C = $null; D = $null; E = "  +2.32%  " }
    @{ Row = 48; B = $null; C = $null; D = "21.06"; E = "  +2.67%  " }
    @{ Row = 49; B = $null; C = $null; D = "0.0531"; E = "  -0.60%  " }
    @{ Row = 50; B = $null; C = $null; D = "0.599"; E = "  -0.29%  " }
    @{ Row = 51; B = $null; C = $null; D = $null; E = "  -0.03%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    # Price (D) and Volume (E) columns hold numeric-looking text (e.g. "7.30",
    # "0.0000101", "  +0.69%  "); force Text format first so Excel stores them
    # as strings instead of coercing/rounding them into numbers.
    if ($u.D -ne $null) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.E -ne $null) {
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $u.E
    }
    if ($u.B -ne $null) { $ws.Range("B$r").Value = $u.B }
    if ($u.C -ne $null) { $ws.Range("C$r").Value = $u.C }
}
